$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Pawan Mankar" -> "manda user" for every Pass row's Message column
$ws.Range("D2").Value = "manda user"
$ws.Range("D5").Value = "manda user"
$ws.Range("D7").Value = "manda user"

# Format the Result column: Pass rows -> bold green, Fail rows -> italic red
foreach ($r in 2,5,7) {
    $ws.Range("C$r").Font.Bold = $true
    $ws.Range("C$r").Font.ColorIndex = 10
}

foreach ($r in 3,4,6,8) {
    $ws.Range("C$r").Font.Italic = $true
    $ws.Range("C$r").Font.ColorIndex = 3
}

# Move the active selection to C3
$ws.Range("C3").Select() | Out-Null
